# Applies the commit "Update master to output generated at c8c62b6":
# - update the date line
# - update the division-practice answers in the table
# Most answer strings are unique in the document, so a straightforward
# Find/Replace is used for them. One value ("37÷8=4, 5") appears twice
# with two different replacements, so those two cells are addressed
# directly via the table/cell index to avoid ambiguity.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Header date
Replace-Text "2025-09-29 Monday" "2025-09-30 Tuesday"

# Row 1 (table row 1)
Replace-Text "27÷6=4, 3"  "71÷9=7, 8"
Replace-Text "41÷2=20, 1" "50÷6=8, 2"
Replace-Text "51÷9=5, 6"  "81÷8=10, 1"
Replace-Text "37÷5=7, 2"  "46÷2=23, 0"
Replace-Text "66÷2=33, 0" "21÷8=2, 5"

# Row 2 (table row 5)
Replace-Text "29÷2=14, 1" "46÷5=9, 1"
Replace-Text "43÷6=7, 1"  "95÷8=11, 7"
Replace-Text "88÷4=22, 0" "82÷4=20, 2"
Replace-Text "82÷6=13, 4" "54÷9=6, 0"
Replace-Text "70÷5=14, 0" "90÷2=45, 0"

# Row 3 (table row 9) - column 3 ("37÷8=4, 5") is ambiguous with row 4's
# column 4, so address both occurrences directly by cell.
Replace-Text "57÷7=8, 1" "25÷5=5, 0"
Replace-Text "16÷6=2, 4" "97÷8=12, 1"
$d.Tables.Item(1).Cell(9, 3).Range.Text = "83÷8=10, 3"
Replace-Text "58÷8=7, 2" "81÷4=20, 1"
Replace-Text "40÷3=13, 1" "49÷8=6, 1"

# Row 4 (table row 13)
Replace-Text "38÷7=5, 3" "87÷8=10, 7"
Replace-Text "80÷4=20, 0" "73÷3=24, 1"
Replace-Text "12÷3=4, 0" "11÷2=5, 1"
$d.Tables.Item(1).Cell(13, 4).Range.Text = "57÷3=19, 0"
Replace-Text "94÷2=47, 0" "90÷7=12, 6"

# Row 5 (table row 17)
Replace-Text "63÷9=7, 0" "86÷6=14, 2"
Replace-Text "65÷3=21, 2" "52÷3=17, 1"
Replace-Text "85÷4=21, 1" "68÷2=34, 0"
Replace-Text "68÷5=13, 3" "93÷5=18, 3"
Replace-Text "58÷3=19, 1" "27÷8=3, 3"
